$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal text, bypassing Excel's
# automatic number/date detection that `.Value` performs on plain strings
# (e.g. "2021-07-08" would otherwise turn into a date serial, and "1"
# would otherwise turn into a number).
#
# We build the literal text via a quoted formula in an unused scratch
# cell, copy it, and paste-special only the *value* into the destination
# - this yields a plain string-typed cell with no numeric/date formatting
# and without leaving any extra cell style behind.
$scratch = $ws.Range("AZ1048576")
function Set-TextValue($rangeAddress, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# ---------------- Row 12 ----------------
$ws.Range("A12").Value = 100247925

Set-TextValue "I12" "1"
$ws.Range("J12").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
Set-TextValue "N12" "fönsterfälla"

$ws.Range("S12").Value = 10

Set-TextValue "AA12" "2021-07-08"

$ws.Range("AF12").ClearContents()

Set-TextValue "AI12" "Björkrik klen grannaturskog, sjönära"
Set-TextValue "AO12" "IBL-fälla mellan två björkhögstubbar med fnösketickor"

Set-TextValue "AW12" "Lars-Olof Grund"
Set-TextValue "AX12" "Lars-Olof Grund, Lars-Ove Wikars, Erland Lindblad"
Set-TextValue "AY12" "Uppföljning av naturvärden i lövrika skogar i Lidsjöbergstrakten, Strömsunds kommun, Jämtlands län"

# ---------------- Row 13 ----------------
$ws.Range("A13").Value = 100247928

Set-TextValue "I13" "1"
$ws.Range("J13").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
Set-TextValue "N13" "fönsterfälla"

$ws.Range("S13").Value = 10

Set-TextValue "AA13" "2021-07-08"

$ws.Range("AF13").ClearContents()

Set-TextValue "AI13" "Björkrik klen grannaturskog, sjönära"
Set-TextValue "AO13" "IBL-fälla mellan två björkhögstubbar med fnösketickor"

Set-TextValue "AW13" "Lars-Olof Grund"
Set-TextValue "AX13" "Lars-Olof Grund, Lars-Ove Wikars, Erland Lindblad"
Set-TextValue "AY13" "Uppföljning av naturvärden i lövrika skogar i Lidsjöbergstrakten, Strömsunds kommun, Jämtlands län"

Write-Output "Applied changes to rows 12 and 13"
